$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    5.492632707754335,
    5.499046050334988,
    5.505544663430387,
    5.512129637034153,
    5.518802072123171,
    5.525563080564389,
    5.532413785007606,
    5.539355318763616,
    5.546388825667322,
    5.553515459925096,
    5.560736385945974,
    5.568052778156225,
    5.575465820796686,
    5.58297670770254,
    5.590586642065027,
    5.598296836174752,
    5.606108511146211,
    5.614022896623263,
    5.622041230465212,
    5.630164758413412,
    5.638394733738127,
    5.646732416865642,
    5.65517907498554,
    5.663735981638284,
    5.672404416283124,
    5.681185663846652,
    5.690081014252241,
    5.699091761930777,
    5.708219205313213,
    5.717464646305524,
    5.726829389746727,
    5.736314742850955,
    5.745922014634292,
    5.755652515327709,
    5.765507555777021,
    5.775488446831528,
    5.785596498722468,
    5.795833020433171,
    5.806199319062531,
    5.81669669918377,
    5.82732646220047,
    5.83808990570218,
    5.848988322821869,
    5.860023001597834,
    5.871195224342611,
    5.882506267021777,
    5.893957398645659,
    5.905549880676941,
    5.917284966457522,
    5.929163900658042,
    5.941187918753549,
    5.953358246529019,
    5.965676099618548,
    5.978142683082127,
    5.990759191023956,
    6.003526806256542,
    6.01644670001473,
    6.029520031723973,
    6.042747948827226,
    6.056131586674848,
    6.069672068482035,
    6.083370505358227,
    6.097227996412982,
    6.111245628942823,
    6.12542447870354,
    6.139765610272376,
    6.154270077504337,
    6.16893892408715,
    6.183773184198784,
    6.198773883271787,
    6.213942038868244,
    6.229278661669187,
    6.244784756581867,
    6.260461323968476,
    6.276309360999182,
    6.292329863132553,
    6.308523825725808,
    6.324892245777241,
    6.341436123802834,
    6.358156465848577,
    6.375054285639805,
    6.392130606868521,
    6.409386465618991,
    6.42682291293181,
    6.444441017506025,
    6.462241868538471,
    6.480226578698939,
    6.498396287239532,
    6.516752163235805,
    6.535295408956881,
    6.554027263361425,
    6.572949005715373,
    6.592061959327292,
    6.611367495396441,
    6.63086703696804,
    6.65056206298996,
    6.670454112464273,
    6.690544788686821,
    6.710835763567212,
    6.731328782021339,
    6.752025666428092,
    6.772928321141214,
    6.794038737047009,
    6.815358996158221,
    6.836891276233836,
    6.858637855414363,
    6.880601116861643,
    6.902783553392107,
    6.925187772091941,
    6.94781649890256,
    6.970672583164477,
    6.993759002107403,
    7.017078865274586,
    7.040635418869039,
    7.064432050009299,
    7.088472290882605,
    7.11275982278319,
    7.137298480023459,
    7.162092253706438,
    7.187145295347301,
    7.212461920332851,
    7.238046611207578,
    7.263904020775439,
    7.290038975006997,
    7.316456475741933,
    7.343161703177263,
    7.370160018132588,
    7.397456964083721,
    7.425058268957111,
    7.452969846677969,
    7.481197798465705,
    7.509748413871093,
    7.538628171550246,
    7.567843739771578,
    7.597401976652303,
    7.627309930122411,
    7.657574837614542,
    7.688204125479361,
    7.719205408126724,
    7.750586486894208,
    7.782355348645056,
    7.814520164099009,
    7.84708928590016,
    7.88007124642705,
    7.913474755351005,
    7.947308696949965,
    7.981582127185651,
    8.016304270552958,
    8.051484516711485,
    8.087132416909618,
    8.123257680212797,
    8.159870169548185,
    8.196979897578714,
    8.234597022420346,
    8.272731843216798,
    8.311394795587077,
    8.350596446961172,
    8.390347491820386,
    8.430658746858727,
    8.471541146082723,
    8.51300573586698,
    8.555063669983459,
    8.597726204622431,
    8.641004693423508,
    8.684910582535139,
    8.729455405721207,
    8.774650779533134,
    8.820508398566165,
    8.867040030818293,
    8.914257513169924,
    8.962172747002628,
    9.010797693974688,
    9.06014437197091,
    9.110224851243929,
    9.161051250763732,
    9.212635734791608,
    9.264990509694307,
    9.318127821013654,
    9.372059950806031,
    9.426799215265948,
    9.482357962646793,
    9.538748571491448,
    9.595983449184562,
    9.654075030837765,
    9.713035778517941,
    9.77287818082827,
    9.833614752850782,
    9.895258036458268,
    9.957820601002792,
    10.02131504438703,
    10.08575399452392,
    10.15115011118931,
    10.21751608827148,
    10.28486465642053,
    10.35320858609996,
    10.42256069104197,
    10.49293383210713,
    10.56434092154869,
    10.63679492768047,
    10.71030887994752,
    10.78489587439712,
    10.8605690795478,
    10.93734174265085,
    11.01522719634712,
    11.09423886570863,
    11.17439027566255,
    11.25569505879208,
    11.33816696350903,
    11.42181986259215,
    11.50666776208496,
    11.59272481054682,
    11.68000530865033,
    11.7685237191184,
    11.85829467699385,
    11.94933300023436,
    12.04165370062559,
    12.13527199500501,
    12.23020331678925,
    12.32646332779764,
    12.42406793036437,
    12.5230332797323,
    12.62337579672113,
    12.72511218066278,
    12.82825942259711,
    12.93283481872097,
    13.03885598408409,
    13.14634086652509,
    13.25530776084115,
    13.36577532318532,
    13.47776258568505,
    13.59128897127624,
    13.70637430874685,
    13.82303884798438,
    13.94130327542167,
    14.06118872967551,
    14.18271681737262,
    14.30590962915762,
    14.43078975587756,
    14.55738030493776,
    14.68570491682324,
    14.81578778178039,
    14.94765365665323,
    15.08132788186794,
    15.21683639856008,
    15.35420576583754,
    15.49346317817277,
    15.63463648291692,
    15.7777541979281,
    15.92284552930556,
    16.06994038922115,
    16.21906941383794,
    16.37026398130613,
    16.52355622982503,
    16.6789790757585,
    16.83656623179138,
    16.99635222511223,
    17.15837241560713,
    17.32266301404761,
    17.4892611002553,
    17.65820464122282,
    17.82953250917067,
    18.00328449951722,
    18.1795013487368,
    18.3582247520806,
    18.53949738113121,
    18.72336290116111,
    18.90986598826274,
    19.09905234621515,
    19.2909687230509,
    19.48566292728339,
    19.68318384375324,
    19.883581449049,
    20.08690682645532,
    20.29321218037878,
    20.50255085019861,
    20.71497732348687,
    20.93054724853942,
    21.14931744615629,
    21.3713459206066,
    21.59669186971027,
    21.82541569396555,
    22.057579004648,
    22.29324463080353,
    22.53247662505456,
    22.77534026813548,
    23.02190207206983,
    23.27222978189872,
    23.52639237586687,
    23.78446006396886,
    24.04650428475621,
    24.31259770030164,
    24.58281418921464,
    24.85722883760006,
    25.13591792784753,
    25.41895892513851,
    25.70643046155428,
    25.99841231766706,
    26.29498540149349,
    26.59623172468946,
    26.90223437586284,
    27.21307749088004,
    27.43036151308731,
    27.65133001966648,
    27.87605060101795,
    28.10459158275748,
    28.33702199022897,
    28.5734115102248,
    28.8138304498221,
    29.0583496922466,
    29.30704064967667,
    29.55997521290402,
    29.81722569777062,
    30.07886478830471,
    30.3449654764843,
    30.61560099856056,
    30.89084476787984,
    31.17077030414898,
    31.4554511590959,
    31.74496083848418,
    32.03937272044999,
    32.338759970137,
    32.64319545061662,
    32.95275163008832,
    33.2675004853699,
    33.58751340169479,
    33.91286106884985,
    34.24361337369761,
    34.57983928914217,
    34.92160675961173,
    35.2689825831461,
    35.62203229019305,
    35.98082001923359,
    36.34540838937311,
    36.71585837005244,
    37.09222914804999,
    37.47457799196513,
    37.86296011438996,
    38.2574285319951,
    38.65803392377526,
    39.06482448771581,
    39.47784579616278,
    39.89714065019533,
    40.32274893331795,
    40.75470746480757,
    41.19304985306663,
    41.63780634935146,
    42.08900370225768,
    42.546665013363,
    43.01080959443735,
    43.48145282664586,
    43.95860602217913,
    44.44227628875657,
    44.93246639745535,
    45.42917465432502,
    45.93239477625172,
    46.44211577153931,
    46.29090020465513,
    46.13816044785924,
    45.98387264639673,
    45.82801237679073,
    45.67055463539555,
    45.51147382662256,
    45.35074375084763,
    45.18833759201593,
    45.02422790496335,
    44.85838660248418,
    44.69078494217912,
    44.52139351312839,
    44.35018222244286,
    44.17712028175693,
    44.00217619373725,
    43.82531773869347,
    43.64651196139093,
    43.46572515817682,
    43.28292286454829,
    43.09806984330449,
    42.91113007344249,
    42.72206673997237,
    42.53084222484623,
    42.33741809921344,
    42.1417551172339,
    41.94381321170066,
    41.74355149174327,
    41.5409282429033,
    41.33590092989404,
    41.12842620237637,
    40.9184599041027,
    40.7059570858004,
    40.49087202218434,
    40.27315823350673,
    40.05276851206729,
    39.82965495412268,
    39.60376899764617,
    39.37506146640016,
    39.14348262079002,
    38.9089822159752,
    38.67150956771381,
    38.43101362641504,
    38.18744305986858,
    37.94074634510888,
    37.69087186985858,
    37.43776804397314,
    37.18138342128564,
    36.92166683221759,
    36.6585675274862,
    36.39203533319336,
    36.12202081753473,
    35.84847546930835,
    35.57135188834218,
    35.29060398788915,
    35.00618720896394,
    34.7180587465128,
    34.42617778722033,
    34.13050575866265,
    33.83100658941751,
    33.52764697963654,
    33.22039668147666,
    32.90922878867406,
    32.5941200344289,
    32.27505109664999,
    31.95200690948932,
    31.62497697997637,
    31.29395570844248,
    30.95894271130819,
    30.61994314469214,
    30.27696802718916,
    29.93003456006166,
    29.57916644298939,
    29.22439418343433,
    28.86575539759663,
    28.50329510086916,
    28.1370659856412,
    27.7671286842584,
    27.39355201491795,
    27.01641320826456,
    26.63579811245706,
    26.25180137449682,
    25.86452659564923,
    25.47408645884841,
    25.08060282605358,
    24.6842068036237,
    24.28503877389448,
    23.88324839127968,
    23.47899454137482,
    23.07244526171677,
    22.66377762304611,
    22.25317757012906,
    21.84083972142153,
    21.426967127098,
    21.01177098522059,
    20.59547031608723
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

Write-Output "done: $($values.Count) values set"